# UFCconnections.xlsx edit:
#  - Rename the "MEGA" sheet/board to "UFC controller" and relabel its
#    "chip"/"mega" column (A) as "board"/"UFC controller".
#  - Fix the letter "d" segment-code definition on the "16-seg dot matrix"
#    sheet (the 0/1 inputs in B13/D13/C14/D14/B16/C16/B17/D17).
#  - Restore the selections that were active on each sheet when the file
#    was last saved.

$wb = $excel.ActiveWorkbook

# --- "MEGA" sheet -> "UFC controller" -----------------------------------
$wsMega = $wb.Worksheets.Item("MEGA")

# Column header (row 1) "chip" -> "board"
$wsMega.Range("A1").Value = "board"

# Every data row (2-50) was labelled with the chip name "mega";
# relabel to the new board name "UFC controller".
$wsMega.Range("A2:A50").Value = "UFC controller"

# Restore the sheet's last-known selection.
$wsMega.Range("A2:A50").Select()

# Rename the sheet/tab itself (keep using the $wsMega reference so the
# earlier lookup by the old name still resolves correctly).
$wsMega.Name = "UFC controller"

# --- "backlight" sheet: just a selection change -------------------------
$wsBacklight = $wb.Worksheets.Item("backlight")
$wsBacklight.Range("A2:A9").Select()

# --- "16-seg dot matrix" sheet -------------------------------------------
$wsSeg = $wb.Worksheets.Item("16-seg dot matrix")

# Corrected bit pattern describing which segments make up the letter "d".
$wsSeg.Range("B13").Value = 0
$wsSeg.Range("D13").Value = 0
$wsSeg.Range("C14").Value = 1
$wsSeg.Range("D14").Value = 0
$wsSeg.Range("B16").Value = 0
$wsSeg.Range("C16").Value = 1
$wsSeg.Range("B17").Value = 0
$wsSeg.Range("D17").Value = 0

# This sheet was the active tab/selection when the workbook was saved, so
# activate it last and restore its selection.
$wsSeg.Activate()
$wsSeg.Range("T29").Select()
